$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Sheets.Item("ALC")
$updates = @{
    "H28" = 762.2174
    "I28" = 718.94446
    "J28" = 918
    "K28" = 718.94446
    "L28" = 918
    "M28" = -233.94446
    "N28" = -1888
    "H40" = 11318.167
    "I40" = 17294.5
    "J40" = 8330
    "K40" = 17294.5
    "L40" = 8330
    "M40" = -17119.5
    "N40" = -8680
    "H62" = 0
    "I62" = 0
    "J62" = 0
    "K62" = 0
    "L62" = 0
    "H63" = 49949.5
    "J63" = 49949.5
    "L63" = 49949.5
    "N63" = -51197.5
    "H65" = 0
    "I65" = 0
    "J65" = 0
    "K65" = 0
    "L65" = 0
    "H66" = 49949.5
    "J66" = 49949.5
    "L66" = 149848.5
    "N66" = -156088.5
    "H69" = 31743.111
    "I69" = 30137.8
    "J69" = 33749.75
    "K69" = 90413.39999999999
    "L69" = 101249.25
    "M69" = -89539.39999999999
    "N69" = -102997.25
    "H72" = 31743.111
    "I72" = 30137.8
    "J72" = 33749.75
    "K72" = 271240.2
    "L72" = 303747.75
    "M72" = -266872.2
    "N72" = -312483.75
    "H74" = 7576.231
    "I74" = 7576.231
    "K74" = 7576.231
    "M74" = -6640.231
    "H75" = 89442.336
    "J75" = 89442.336
    "L75" = 89442.336
    "N75" = -91314.336
    "H77" = 7576.231
    "I77" = 7576.231
    "K77" = 37881.155
    "M77" = -33201.155
    "H78" = 89442.336
    "J78" = 89442.336
    "L78" = 268327.008
    "N78" = -277687.008
    "H81" = 238558.5
    "I81" = 0
    "K81" = 0
    "H84" = 238558.5
    "I84" = 0
    "K84" = 0
    "H88" = 5900.3716
    "I88" = 2953.3333
    "J88" = 8110.65
    "K88" = 2953.3333
    "L88" = 8110.65
    "M88" = -2547.3333
    "N88" = -8922.65
    "H91" = 5900.3716
    "I91" = 2953.3333
    "J91" = 8110.65
    "K91" = 2953.3333
    "L91" = 8110.65
    "M91" = -1549.3333
    "N91" = -10918.65
    "H96" = 737
    "I96" = 449.6154
    "J96" = 1982.3334
    "K96" = 1348.8462
    "L96" = 5947.0002
    "M96" = 24.15380000000005
    "N96" = -8693.0002
    "H132" = 19361.484
    "I132" = 3025.9375
    "K132" = 9077.8125
    "M132" = -6547.8125
    "H137" = 4073
    "J137" = 4930.3
    "L137" = 14790.9
    "N137" = -19890.9
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
$clears = @("M62", "N62", "M65", "N65", "M81", "M84")
foreach ($key in $clears) {
    $ws.Range($key).ClearContents()
}

# ---- ARM ----
$ws = $wb.Sheets.Item("ARM")
$updates = @{
    "H132" = 2806.4358
    "I132" = 1610.4857
    "J132" = 13271
    "K132" = 4831.4571
    "L132" = 39813
    "M132" = -2301.4571
    "N132" = -44873
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- BSM ----
$ws = $wb.Sheets.Item("BSM")
$updates = @{
    "H105" = 3125.0625
    "I105" = 3235.9285
    "J105" = 2349
    "K105" = 3235.9285
    "L105" = 2349
    "M105" = -1488.9285
    "N105" = -5843
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- CRP ----
$ws = $wb.Sheets.Item("CRP")
$updates = @{
    "H22" = 1413.25
    "I22" = 959.0769
    "J22" = 2256.7144
    "K22" = 959.0769
    "L22" = 2256.7144
    "M22" = -609.0769
    "N22" = -2956.7144
    "H141" = 140453.67
    "J141" = 140453.67
    "L141" = 140453.67
    "N141" = -150813.67
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- CUL ----
$ws = $wb.Sheets.Item("CUL")
$updates = @{
    "H12" = 56.304348
    "J12" = 50.125
    "L12" = 150.375
    "N12" = -496.375
    "H68" = 1989
    "J68" = 1989
    "L68" = 5967
    "N68" = -7589
    "H71" = 1989
    "J71" = 1989
    "L71" = 17901
    "N71" = -26013
    "H107" = 1451.5
    "I107" = 903
    "J107" = 2000
    "K107" = 2709
    "L107" = 6000
    "M107" = -789
    "N107" = -9840
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- GSM ----
$ws = $wb.Sheets.Item("GSM")
$updates = @{
    "H122" = 5000
    "I122" = 5000
    "K122" = 15000
    "M122" = -12550
    "H128" = 0
    "J128" = 0
    "L128" = 0
    "H132" = 4358.143
    "I132" = 4340.6606
    "J132" = 4428.0713
    "K132" = 13021.9818
    "L132" = 13284.2139
    "M132" = -10491.9818
    "N132" = -18344.2139
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
$clears = @("N128")
foreach ($key in $clears) {
    $ws.Range($key).ClearContents()
}

# ---- LTW ----
$ws = $wb.Sheets.Item("LTW")
$updates = @{
    "H46" = 1880.3422
    "J46" = 2177.6072
    "L46" = 2177.6072
    "N46" = -2553.6072
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- WVR ----
$ws = $wb.Sheets.Item("WVR")
$updates = @{
    "H41" = 9428.75
    "I41" = 9670
    "J41" = 9187.5
    "K41" = 9670
    "L41" = 9187.5
    "M41" = -9280
    "N41" = -9967.5
    "H132" = 4165.0986
    "I132" = 3456.5085
    "K132" = 10369.5255
    "M132" = -7839.5255
    "H136" = 1673.4688
    "I136" = 1606.2273
    "K136" = 4818.6819
    "M136" = -2268.6819
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
